$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The Price/Volume columns store plain numeric- and percent-looking
# strings as text (inline strings in the source workbook), so each
# write forces text format, assigns the literal value, then clears
# the temporary formatting so the cell keeps its original (default)
# style while remaining text-typed.
function Set-TextValue($cell, $val) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "300.62"
Set-TextValue "E2" "-0.95%"
Set-TextValue "D3" "31.42"
Set-TextValue "E3" "-1.97%"
Set-TextValue "E4" "-2.33%"
Set-TextValue "D5" "0.07373"
Set-TextValue "E5" "-1.43%"
Set-TextValue "D6" "2.501"
Set-TextValue "E6" "64.25%"
Set-TextValue "D7" "7.937"
Set-TextValue "E7" "1.12%"
Set-TextValue "D8" "3.767"
Set-TextValue "E8" "-1.09%"
Set-TextValue "D9" "0.9227"
Set-TextValue "E9" "0.48%"
Set-TextValue "D10" "0.1727"
Set-TextValue "E10" "2.59%"
Set-TextValue "D11" "0.07599"
Set-TextValue "E11" "-5.38%"
Set-TextValue "D12" "0.08067"
Set-TextValue "E12" "0.03%"
Set-TextValue "D13" "0.03032"
Set-TextValue "E13" "0.64%"
Set-TextValue "D14" "0.09926"
Set-TextValue "E14" "0.21%"
Set-TextValue "D15" "0.001492"
Set-TextValue "E15" "0.02%"
Set-TextValue "D16" "0.006083"
Set-TextValue "E16" "-4.59%"
Set-TextValue "E17" "-0.25%"
Set-TextValue "E18" "-0.10%"
Set-TextValue "E19" "-0.69%"
Set-TextValue "D20" "0.1336"
Set-TextValue "E20" "0.02%"
Set-TextValue "D21" "4.657"
Set-TextValue "E21" "4.06%"
Set-TextValue "D22" "0.04651"
Set-TextValue "E22" "1.17%"
Set-TextValue "D23" "0.1566"
Set-TextValue "E23" "-3.29%"
Set-TextValue "D24" "0.001223"
Set-TextValue "E24" "0.47%"
Set-TextValue "D25" "0.004492"
Set-TextValue "E25" "1.05%"
Set-TextValue "D26" "0.0001300"
Set-TextValue "E26" "-7.05%"
Set-TextValue "E27" "5.46%"
Set-TextValue "D39" "0.01731"
Set-TextValue "E39" "0.88%"
Set-TextValue "D40" "0.04523"
Set-TextValue "E40" "0.37%"
Set-TextValue "D41" "0.007172"
Set-TextValue "E41" "0.43%"
Set-TextValue "D42" "0.1349"
Set-TextValue "E42" "0.12%"
Set-TextValue "D43" "0.002229"
Set-TextValue "E43" "-0.79%"
Set-TextValue "D44" "0.01073"
Set-TextValue "E44" "-16.41%"
Set-TextValue "D45" "0.00006279"
Set-TextValue "E45" "1.72%"
Set-TextValue "D46" "0.8083"
Set-TextValue "E46" "13.89%"
Set-TextValue "E47" "-22.95%"
